$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.233.22"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "3.865.91"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'690.66"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("D6").Value = "'173.20"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").Value = "3.862.72"
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").Value = "'7.40"
$ws.Range("E11").Value = "  +4.50%  "
$ws.Range("D12").Value = "'0.462"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  +6.05%  "
$ws.Range("D14").Value = "'36.65"
$ws.Range("E14").Value = "  +2.82%  "
$ws.Range("D15").Value = "4.505.99"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "3.859.34"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "71.272.87"
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("D18").Value = "'17.81"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "'7.25"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'11.09"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "'490.33"
$ws.Range("E22").Value = "  +3.49%  "
$ws.Range("D23").Value = "'0.723"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").Value = "'84.85"
$ws.Range("E24").Value = "  +2.26%  "
$ws.Range("D25").Value = "'0.0000148"
$ws.Range("E25").Value = "  +4.44%  "
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("E27").Value = "  +2.05%  "
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "4.014.05"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'3.12"
$ws.Range("E31").Value = "  +9.55%  "
$ws.Range("D32").Value = "'7.65"
$ws.Range("E32").Value = "  +2.79%  "
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "'29.85"
$ws.Range("E34").Value = "  +1.18%  "
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").Value = "'9.31"
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("D37").Value = "3.814.29"
$ws.Range("E37").Value = "  +1.75%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("E40").Value = "  +12.88%  "
$ws.Range("D41").Value = "'3.45"
$ws.Range("E41").Value = "  +1.70%  "
$ws.Range("D42").Value = "'6.08"
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("E43").Value = "  +6.45%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'163.77"
$ws.Range("E46").Value = "  +3.06%  "
$ws.Range("E47").Value = "  +7.76%  "
$ws.Range("D48").Value = "'48.73"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "'44.51"
$ws.Range("E49").Value = "  -2.23%  "
$ws.Range("D50").Value = "'0.304"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("E51").Value = "  -1.94%  "
